# Rename sheets (shift names forward: Sheet1 -> Sheet2, Sheet2 -> Sheet3)
# and update the selected/active cell on the first sheet from B3 to C3.

$wb = $excel.ActiveWorkbook

$sheet1 = $wb.Worksheets.Item(1)
$sheet2 = $wb.Worksheets.Item(2)

# Rename the second sheet first to avoid a name collision while renaming.
$sheet2.Name = "Sheet3"
$sheet1.Name = "Sheet2"

# Update the selection on the (renamed) first sheet from B3 to C3,
# then restore the originally active sheet so the workbook's active
# tab is unaffected by this housekeeping change.
$originallyActive = $wb.ActiveSheet
$sheet1.Activate()
$sheet1.Range("C3").Select()
$originallyActive.Activate()
